$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column A for "Date" (shifts Product -> B, Sales -> C) ---
$ws.Columns("A").Insert()

# --- Copy the header style (bold font + border + centered alignment) from the
#     existing "Product" header (now in B1) onto the new header cells ---
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("E1").PasteSpecial(-4122)

# --- Header row ---
$ws.Range("A1").Value = "Date"
$ws.Range("D1").Value = "Units of Sale"
$ws.Range("E1").Value = "Cost per Unit"

# Register numFmtId 164 ("yyyy-mm-dd h:mm:ss") by applying it to a single
# cell first, then apply the final, differently-cased format code
# ("YYYY-MM-DD HH:MM:SS", numFmtId 165) to the whole date range - this is
# the format that actually ends up on the date cells.
$ws.Range("A2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A2:A9").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- Date column values (serial date numbers) ---
$ws.Range("A2").Value = 45974
$ws.Range("A3").Value = 45978
$ws.Range("A4").Value = 45982
$ws.Range("A5").Value = 45985
$ws.Range("A6").Value = 46033
$ws.Range("A7").Value = 46040
$ws.Range("A8").Value = 46048
$ws.Range("A9").Value = 46053

# --- Units of Sale column ---
$ws.Range("D2").Value = 5
$ws.Range("D3").Value = 25
$ws.Range("D4").Value = 18
$ws.Range("D5").Value = 12
$ws.Range("D6").Value = 45
$ws.Range("D7").Value = 50
$ws.Range("D8").Value = 32
$ws.Range("D9").Value = 28

# --- Cost per Unit column ---
$ws.Range("E2").Value = 9000
$ws.Range("E3").Value = 2680
$ws.Range("E4").Value = 1777.78
$ws.Range("E5").Value = 2333.33
$ws.Range("E6").Value = 333.33
$ws.Range("E7").Value = 240
$ws.Range("E8").Value = 718.75
$ws.Range("E9").Value = 642.86
